# Insert a new data row before the current row 128 ("Hortaliza, Vega Modelo de
# Temuco - Ciboulette" weekly update). Excel shifts all the rows from 128..155
# down to 129..156 and copies the formatting (including the date style on
# column D) from the row above, which matches the target workbook layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(128).Insert()

$ws.Cells.Item(128, 1).Value  = 10
$ws.Cells.Item(128, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(128, 3).Value  = 'La Araucanía'
$ws.Cells.Item(128, 4).Value  = 44476
$ws.Cells.Item(128, 5).Value  = 9
$ws.Cells.Item(128, 6).Value  = 100112039
$ws.Cells.Item(128, 7).Value  = 'Ciboulette'
$ws.Cells.Item(128, 8).Value  = 'Sin especificar'
$ws.Cells.Item(128, 9).Value  = 'Primera'
$ws.Cells.Item(128, 10).Value = 20
$ws.Cells.Item(128, 11).Value = 7000
$ws.Cells.Item(128, 12).Value = 7000
$ws.Cells.Item(128, 13).Value = 7000
$ws.Cells.Item(128, 14).Value = '$/docena de atados'
$ws.Cells.Item(128, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(128, 16).Value = 2333
$ws.Cells.Item(128, 17).Value = 3
$ws.Cells.Item(128, 18).Value = 'Hortaliza'
